$d = $word.ActiveDocument

# Paragraphs (1-indexed via Word COM) whose heading text needs to become bold:
#   1  - "Nestor Wilke"
#   5  - "Experiência de trabalho"
#   6  - "Gerente de equipe de animação"
#   12 - "Designer sênior de animação" (also renamed below)
#   18 - "Designer de animação"
#   24 - "Bacharelado em Belas Artes em Animação"
$boldParagraphIndexes = @(1, 5, 6, 12, 18, 24)

foreach ($idx in $boldParagraphIndexes) {
    $p = $d.Paragraphs($idx)
    $p.Range.Font.Bold = 1
}

# Reorder the job-title heading text (paragraph 12), scoped so the other two
# occurrences of the phrase elsewhere in the document stay untouched.
$titleRange = $d.Paragraphs(12).Range
$titleRange.Find.Execute("Designer sênior de animação", $true, $true, $false, $false, $false, `
                          $true, 1, $false, "Designer de animação sênior", 2)
